# =========================================================
# Daily Report update: 2026-02-12 (adds 2026-02-11 daily rows)
# =========================================================

$wb = $excel.ActiveWorkbook
$wsDaily = $wb.Worksheets.Item("Daily_Data")
$wsSummary = $wb.Worksheets.Item("Today_Summary")
$wsMonthly = $wb.Worksheets.Item("Monthly_Stats")

# ---- Sheet1 (Daily_Data): append new daily rows 74-97 ----

$wsDaily.Cells.Item(74,1).Value = 46064
$wsDaily.Cells.Item(74,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(74,2).Value = 'ASAHI DEPOSITORY LLC Registered'
$wsDaily.Cells.Item(74,3).Value = 23953631.592
$wsDaily.Cells.Item(74,4).Value = 0
$wsDaily.Cells.Item(74,5).Value = 0
$wsDaily.Cells.Item(74,6).Value = 0
$wsDaily.Cells.Item(74,7).Value = 0
$wsDaily.Cells.Item(74,8).Value = 23953631.592

$wsDaily.Cells.Item(75,1).Value = 46064
$wsDaily.Cells.Item(75,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(75,2).Value = 'ASAHI DEPOSITORY LLC Eligible'
$wsDaily.Cells.Item(75,3).Value = 2987818.038
$wsDaily.Cells.Item(75,4).Value = 0
$wsDaily.Cells.Item(75,5).Value = 431920.43
$wsDaily.Cells.Item(75,6).Value = -431920.43
$wsDaily.Cells.Item(75,7).Value = 0
$wsDaily.Cells.Item(75,8).Value = 2555897.608

$wsDaily.Cells.Item(76,1).Value = 46064
$wsDaily.Cells.Item(76,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(76,2).Value = 'BRINK''S, INC. Registered'
$wsDaily.Cells.Item(76,3).Value = 17562594.449
$wsDaily.Cells.Item(76,4).Value = 0
$wsDaily.Cells.Item(76,5).Value = 0
$wsDaily.Cells.Item(76,6).Value = 0
$wsDaily.Cells.Item(76,7).Value = -1440234.803
$wsDaily.Cells.Item(76,8).Value = 16122359.646

$wsDaily.Cells.Item(77,1).Value = 46064
$wsDaily.Cells.Item(77,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(77,2).Value = 'BRINK''S, INC. Eligible'
$wsDaily.Cells.Item(77,3).Value = 39199825.671
$wsDaily.Cells.Item(77,4).Value = 0
$wsDaily.Cells.Item(77,5).Value = 0
$wsDaily.Cells.Item(77,6).Value = 0
$wsDaily.Cells.Item(77,7).Value = 1440234.803
$wsDaily.Cells.Item(77,8).Value = 40640060.474

$wsDaily.Cells.Item(78,1).Value = 46064
$wsDaily.Cells.Item(78,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(78,2).Value = 'CNT DEPOSITORY, INC. Registered'
$wsDaily.Cells.Item(78,3).Value = 15343500.119
$wsDaily.Cells.Item(78,4).Value = 0
$wsDaily.Cells.Item(78,5).Value = 0
$wsDaily.Cells.Item(78,6).Value = 0
$wsDaily.Cells.Item(78,7).Value = -2368902.04
$wsDaily.Cells.Item(78,8).Value = 12974598.079

$wsDaily.Cells.Item(79,1).Value = 46064
$wsDaily.Cells.Item(79,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(79,2).Value = 'CNT DEPOSITORY, INC. Eligible'
$wsDaily.Cells.Item(79,3).Value = 12937863.863
$wsDaily.Cells.Item(79,4).Value = 0
$wsDaily.Cells.Item(79,5).Value = 0
$wsDaily.Cells.Item(79,6).Value = 0
$wsDaily.Cells.Item(79,7).Value = 2368902.04
$wsDaily.Cells.Item(79,8).Value = 15306765.903

$wsDaily.Cells.Item(80,1).Value = 46064
$wsDaily.Cells.Item(80,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(80,2).Value = 'DELAWARE DEPOSITORY Registered'
$wsDaily.Cells.Item(80,3).Value = 1966294.501
$wsDaily.Cells.Item(80,4).Value = 0
$wsDaily.Cells.Item(80,5).Value = 0
$wsDaily.Cells.Item(80,6).Value = 0
$wsDaily.Cells.Item(80,7).Value = -413592.568
$wsDaily.Cells.Item(80,8).Value = 1552701.933

$wsDaily.Cells.Item(81,1).Value = 46064
$wsDaily.Cells.Item(81,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(81,2).Value = 'DELAWARE DEPOSITORY Eligible'
$wsDaily.Cells.Item(81,3).Value = 15977369.624
$wsDaily.Cells.Item(81,4).Value = 0
$wsDaily.Cells.Item(81,5).Value = 129734.436
$wsDaily.Cells.Item(81,6).Value = -129734.436
$wsDaily.Cells.Item(81,7).Value = 413592.568
$wsDaily.Cells.Item(81,8).Value = 16261227.756

$wsDaily.Cells.Item(82,1).Value = 46064
$wsDaily.Cells.Item(82,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(82,2).Value = 'HSBC BANK, USA Registered'
$wsDaily.Cells.Item(82,3).Value = 3472271.68
$wsDaily.Cells.Item(82,4).Value = 0
$wsDaily.Cells.Item(82,5).Value = 0
$wsDaily.Cells.Item(82,6).Value = 0
$wsDaily.Cells.Item(82,7).Value = 0
$wsDaily.Cells.Item(82,8).Value = 3472271.68

$wsDaily.Cells.Item(83,1).Value = 46064
$wsDaily.Cells.Item(83,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(83,2).Value = 'HSBC BANK, USA Eligible'
$wsDaily.Cells.Item(83,3).Value = 21150312.483
$wsDaily.Cells.Item(83,4).Value = 0
$wsDaily.Cells.Item(83,5).Value = 0
$wsDaily.Cells.Item(83,6).Value = 0
$wsDaily.Cells.Item(83,7).Value = 0
$wsDaily.Cells.Item(83,8).Value = 21150312.483

$wsDaily.Cells.Item(84,1).Value = 46064
$wsDaily.Cells.Item(84,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(84,2).Value = 'INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered'
$wsDaily.Cells.Item(84,3).Value = 620749.47
$wsDaily.Cells.Item(84,4).Value = 0
$wsDaily.Cells.Item(84,5).Value = 0
$wsDaily.Cells.Item(84,6).Value = 0
$wsDaily.Cells.Item(84,7).Value = -346959.6
$wsDaily.Cells.Item(84,8).Value = 273789.87

$wsDaily.Cells.Item(85,1).Value = 46064
$wsDaily.Cells.Item(85,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(85,2).Value = 'INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible'
$wsDaily.Cells.Item(85,3).Value = 3295246.644
$wsDaily.Cells.Item(85,4).Value = 0
$wsDaily.Cells.Item(85,5).Value = 0
$wsDaily.Cells.Item(85,6).Value = 0
$wsDaily.Cells.Item(85,7).Value = 346959.6
$wsDaily.Cells.Item(85,8).Value = 3642206.244

$wsDaily.Cells.Item(86,1).Value = 46064
$wsDaily.Cells.Item(86,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(86,2).Value = 'JP MORGAN CHASE BANK NA Registered'
$wsDaily.Cells.Item(86,3).Value = 12035257.32
$wsDaily.Cells.Item(86,4).Value = 0
$wsDaily.Cells.Item(86,5).Value = 0
$wsDaily.Cells.Item(86,6).Value = 0
$wsDaily.Cells.Item(86,7).Value = 0
$wsDaily.Cells.Item(86,8).Value = 12035257.32

$wsDaily.Cells.Item(87,1).Value = 46064
$wsDaily.Cells.Item(87,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(87,2).Value = 'JP MORGAN CHASE BANK NA Eligible'
$wsDaily.Cells.Item(87,3).Value = 151278543.683
$wsDaily.Cells.Item(87,4).Value = 0
$wsDaily.Cells.Item(87,5).Value = 1126603.4
$wsDaily.Cells.Item(87,6).Value = -1126603.4
$wsDaily.Cells.Item(87,7).Value = 0
$wsDaily.Cells.Item(87,8).Value = 150151940.283

$wsDaily.Cells.Item(88,1).Value = 46064
$wsDaily.Cells.Item(88,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(88,2).Value = 'LOOMIS INTERNATIONAL (US) LLC Registered'
$wsDaily.Cells.Item(88,3).Value = 7374299.767
$wsDaily.Cells.Item(88,4).Value = 0
$wsDaily.Cells.Item(88,5).Value = 0
$wsDaily.Cells.Item(88,6).Value = 0
$wsDaily.Cells.Item(88,7).Value = 0
$wsDaily.Cells.Item(88,8).Value = 7374299.767

$wsDaily.Cells.Item(89,1).Value = 46064
$wsDaily.Cells.Item(89,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(89,2).Value = 'LOOMIS INTERNATIONAL (US) LLC Eligible'
$wsDaily.Cells.Item(89,3).Value = 23345860.186
$wsDaily.Cells.Item(89,4).Value = 0
$wsDaily.Cells.Item(89,5).Value = 50476.75
$wsDaily.Cells.Item(89,6).Value = -50476.75
$wsDaily.Cells.Item(89,7).Value = 0
$wsDaily.Cells.Item(89,8).Value = 23295383.436

$wsDaily.Cells.Item(90,1).Value = 46064
$wsDaily.Cells.Item(90,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(90,2).Value = 'MALCA-AMIT ARMORED, INC. Registered'
$wsDaily.Cells.Item(90,3).Value = 0
$wsDaily.Cells.Item(90,4).Value = 0
$wsDaily.Cells.Item(90,5).Value = 0
$wsDaily.Cells.Item(90,6).Value = 0
$wsDaily.Cells.Item(90,7).Value = 0
$wsDaily.Cells.Item(90,8).Value = 0

$wsDaily.Cells.Item(91,1).Value = 46064
$wsDaily.Cells.Item(91,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(91,2).Value = 'MALCA-AMIT ARMORED, INC. Eligible'
$wsDaily.Cells.Item(91,3).Value = 0
$wsDaily.Cells.Item(91,4).Value = 0
$wsDaily.Cells.Item(91,5).Value = 0
$wsDaily.Cells.Item(91,6).Value = 0
$wsDaily.Cells.Item(91,7).Value = 0
$wsDaily.Cells.Item(91,8).Value = 0

$wsDaily.Cells.Item(92,1).Value = 46064
$wsDaily.Cells.Item(92,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(92,2).Value = 'MALCA-AMIT USA, LLC Registered'
$wsDaily.Cells.Item(92,3).Value = 1225506.264
$wsDaily.Cells.Item(92,4).Value = 0
$wsDaily.Cells.Item(92,5).Value = 0
$wsDaily.Cells.Item(92,6).Value = 0
$wsDaily.Cells.Item(92,7).Value = 0
$wsDaily.Cells.Item(92,8).Value = 1225506.264

$wsDaily.Cells.Item(93,1).Value = 46064
$wsDaily.Cells.Item(93,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(93,2).Value = 'MALCA-AMIT USA, LLC Eligible'
$wsDaily.Cells.Item(93,3).Value = 798026.177
$wsDaily.Cells.Item(93,4).Value = 0
$wsDaily.Cells.Item(93,5).Value = 0
$wsDaily.Cells.Item(93,6).Value = 0
$wsDaily.Cells.Item(93,7).Value = 0
$wsDaily.Cells.Item(93,8).Value = 798026.177

$wsDaily.Cells.Item(94,1).Value = 46064
$wsDaily.Cells.Item(94,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(94,2).Value = 'MANFRA, TORDELLA & BROOKES, LLC Registered'
$wsDaily.Cells.Item(94,3).Value = 7038609.2
$wsDaily.Cells.Item(94,4).Value = 0
$wsDaily.Cells.Item(94,5).Value = 0
$wsDaily.Cells.Item(94,6).Value = 0
$wsDaily.Cells.Item(94,7).Value = -538131.579
$wsDaily.Cells.Item(94,8).Value = 6500477.621

$wsDaily.Cells.Item(95,1).Value = 46064
$wsDaily.Cells.Item(95,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(95,2).Value = 'MANFRA, TORDELLA & BROOKES, LLC Eligible'
$wsDaily.Cells.Item(95,3).Value = 12226734.004
$wsDaily.Cells.Item(95,4).Value = 0
$wsDaily.Cells.Item(95,5).Value = 597061.864
$wsDaily.Cells.Item(95,6).Value = -597061.864
$wsDaily.Cells.Item(95,7).Value = 538131.579
$wsDaily.Cells.Item(95,8).Value = 12167803.719

$wsDaily.Cells.Item(96,1).Value = 46064
$wsDaily.Cells.Item(96,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(96,2).Value = 'STONEX PRECIOUS METALS LLC Registered'
$wsDaily.Cells.Item(96,3).Value = 7545291.14
$wsDaily.Cells.Item(96,4).Value = 0
$wsDaily.Cells.Item(96,5).Value = 0
$wsDaily.Cells.Item(96,6).Value = 0
$wsDaily.Cells.Item(96,7).Value = 0
$wsDaily.Cells.Item(96,8).Value = 7545291.14

$wsDaily.Cells.Item(97,1).Value = 46064
$wsDaily.Cells.Item(97,1).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$wsDaily.Cells.Item(97,2).Value = 'STONEX PRECIOUS METALS LLC Eligible'
$wsDaily.Cells.Item(97,3).Value = 233197.38
$wsDaily.Cells.Item(97,4).Value = 0
$wsDaily.Cells.Item(97,5).Value = 0
$wsDaily.Cells.Item(97,6).Value = 0
$wsDaily.Cells.Item(97,7).Value = 0
$wsDaily.Cells.Item(97,8).Value = 233197.38
# ---- Sheet2 (Today_Summary): refresh Eligible/Registered/Total_Stock ----

$wsSummary.Cells.Item(2,2).Value = 2555897.608
$wsSummary.Cells.Item(2,4).Value = 26509529.2

$wsSummary.Cells.Item(3,2).Value = 40640060.474
$wsSummary.Cells.Item(3,3).Value = 16122359.646

$wsSummary.Cells.Item(4,2).Value = 15306765.903
$wsSummary.Cells.Item(4,3).Value = 12974598.079

$wsSummary.Cells.Item(5,2).Value = 16261227.756
$wsSummary.Cells.Item(5,3).Value = 1552701.933
$wsSummary.Cells.Item(5,4).Value = 17813929.689

$wsSummary.Cells.Item(7,2).Value = 3642206.244
$wsSummary.Cells.Item(7,3).Value = 273789.87

$wsSummary.Cells.Item(8,2).Value = 150151940.283
$wsSummary.Cells.Item(8,4).Value = 162187197.603

$wsSummary.Cells.Item(9,2).Value = 23295383.436
$wsSummary.Cells.Item(9,4).Value = 30669683.203

$wsSummary.Cells.Item(12,2).Value = 12167803.719
$wsSummary.Cells.Item(12,3).Value = 6500477.621
$wsSummary.Cells.Item(12,4).Value = 18668281.34
# ---- Sheet3 (Monthly_Stats): refresh monthly aggregates ----

# Grand total row (row 2)
$wsMonthly.Cells.Item(2,2).Value = 286202821.463
$wsMonthly.Cells.Item(2,3).Value = 93030184.912
$wsMonthly.Cells.Item(2,4).Value = 379233006.375

# Per-depository detail rows (7-28)
$wsMonthly.Cells.Item(7,4).Value = 1671531.13
$wsMonthly.Cells.Item(7,5).Value = 2555897.608

$wsMonthly.Cells.Item(9,5).Value = 40640060.474

$wsMonthly.Cells.Item(10,5).Value = 16122359.646

$wsMonthly.Cells.Item(11,5).Value = 15306765.903

$wsMonthly.Cells.Item(12,5).Value = 12974598.079

$wsMonthly.Cells.Item(13,4).Value = 174999.025
$wsMonthly.Cells.Item(13,5).Value = 16261227.756

$wsMonthly.Cells.Item(14,5).Value = 1552701.933

$wsMonthly.Cells.Item(17,5).Value = 3642206.244

$wsMonthly.Cells.Item(18,5).Value = 273789.87

$wsMonthly.Cells.Item(19,4).Value = 6676947.6
$wsMonthly.Cells.Item(19,5).Value = 150151940.283

$wsMonthly.Cells.Item(21,4).Value = 2640249.5
$wsMonthly.Cells.Item(21,5).Value = 23295383.436

$wsMonthly.Cells.Item(27,4).Value = 849925.8229999999
$wsMonthly.Cells.Item(27,5).Value = 12167803.719

$wsMonthly.Cells.Item(28,5).Value = 6500477.621
